# Applies the edit described in the commit:
# - Adds a new shared string "Código para lectura/escritura del SEL"
# - Fills in row 18 of the "Metricas" sheet with data for that new task
# - Moves the active cell selection to C24
# All dependent formulas (G18, J18, B24:J24, B25, B28:C33, B34, chart cache)
# recalculate automatically from the new input values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metricas")

# New task row (row 18)
$ws.Range("A18").Value = "Código para lectura/escritura del SEL"
$ws.Range("B18").Value = 70
$ws.Range("C18").Value = 86
$ws.Range("D18").Value = 0.013888888888888888
$ws.Range("E18").Value = 0.1875
$ws.Range("F18").Value = 0.21666666666666667
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0

# Force full recalculation so dependent formulas / chart caches update
$excel.CalculateFullRebuild()

# Update the selected/active cell shown when the workbook is opened
$ws.Activate()
$ws.Range("C24").Select()
